$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
# row 18
$ws.Range("H18").Value = 5923.278
$ws.Range("I18").Value = 1430.5454
$ws.Range("K18").Value = 1430.5454
$ws.Range("M18").Value = -1146.5454
# row 58
$ws.Range("H58").Value = 730.4
$ws.Range("I58").Value = 613
$ws.Range("K58").Value = 1839
$ws.Range("M58").Value = -1689
# row 62
$ws.Range("H62").Value = 1700
$ws.Range("I62").Value = 1700
$ws.Range("K62").Value = 1700
$ws.Range("M62").Value = -1076
# row 65
$ws.Range("H65").Value = 1700
$ws.Range("I65").Value = 1700
$ws.Range("K65").Value = 8500
$ws.Range("M65").Value = -5380
# row 76
$ws.Range("H76").Value = 3906487.2
$ws.Range("I76").Value = 5858106
$ws.Range("J76").Value = 3250
$ws.Range("K76").Value = 5858106
$ws.Range("L76").Value = 3250
$ws.Range("M76").Value = -5857791
$ws.Range("N76").Value = -3880
# row 79
$ws.Range("H79").Value = 3906487.2
$ws.Range("I79").Value = 5858106
$ws.Range("J79").Value = 3250
$ws.Range("K79").Value = 5858106
$ws.Range("L79").Value = 3250
$ws.Range("M79").Value = -5857014
$ws.Range("N79").Value = -5434
# row 95
$ws.Range("H95").Value = 41523.6
$ws.Range("J95").Value = 41523.6
$ws.Range("L95").Value = 41523.6
$ws.Range("N95").Value = -47015.6
# row 116
$ws.Range("H116").Value = 18355.555
$ws.Range("I116").Value = 51350
$ws.Range("J116").Value = 8928.571
$ws.Range("K116").Value = 51350
$ws.Range("L116").Value = 8928.571
$ws.Range("M116").Value = -47908
$ws.Range("N116").Value = -15812.571
# row 132
$ws.Range("H132").Value = 1124.7593
$ws.Range("I132").Value = 1022.95746
$ws.Range("K132").Value = 3068.87238
$ws.Range("M132").Value = -538.8723799999998
# row 137
$ws.Range("H137").Value = 1163.4073
$ws.Range("I137").Value = 790.2941
$ws.Range("K137").Value = 2370.8823
$ws.Range("M137").Value = 179.1177000000002
# row 141
$ws.Range("H141").Value = 701332.8
$ws.Range("J141").Value = 4676.3335
$ws.Range("L141").Value = 14029.0005
$ws.Range("N141").Value = -24389.0005

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
# row 32
$ws.Range("H32").Value = 3788.1528
$ws.Range("I32").Value = 3093.7812
$ws.Range("K32").Value = 3093.7812
$ws.Range("M32").Value = -2806.7812
# row 45
$ws.Range("H45").Value = 1869.2941
$ws.Range("I45").Value = 1932
$ws.Range("K45").Value = 1932
$ws.Range("M45").Value = -1555
# row 61
$ws.Range("H61").Value = 1264.6471
$ws.Range("I61").Value = 628.9677
$ws.Range("K61").Value = 628.9677
$ws.Range("M61").Value = -416.9677
# row 81
$ws.Range("H81").Value = 10000
$ws.Range("I81").Value = 10000
$ws.Range("K81").Value = 10000
$ws.Range("M81").Value = -9002
# row 84
$ws.Range("H84").Value = 10000
$ws.Range("I84").Value = 10000
$ws.Range("K84").Value = 30000
$ws.Range("M84").Value = -25008
# row 101
$ws.Range("H101").Value = 55399.445
$ws.Range("J101").Value = 55399.445
$ws.Range("L101").Value = 55399.445
$ws.Range("N101").Value = -61889.445
# row 122
$ws.Range("H122").Value = 1631.2632
$ws.Range("I122").Value = 1527.875
$ws.Range("J122").Value = 2182.6667
$ws.Range("K122").Value = 4583.625
$ws.Range("L122").Value = 6548.000100000001
$ws.Range("M122").Value = -2133.625
$ws.Range("N122").Value = -11448.0001
# row 132
$ws.Range("H132").Value = 1559.6097
$ws.Range("I132").Value = 1123.8966
$ws.Range("K132").Value = 3371.6898
$ws.Range("M132").Value = -841.6898000000001
# row 136
$ws.Range("H136").Value = 1264.6471
$ws.Range("I136").Value = 628.9677
$ws.Range("K136").Value = 1886.9031
$ws.Range("M136").Value = 663.0969

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
# row 20
$ws.Range("H20").Value = 2668.7856
$ws.Range("I20").Value = 2683.625
$ws.Range("K20").Value = 2683.625
$ws.Range("M20").Value = -2436.625
# row 102
$ws.Range("H102").Value = 10000
$ws.Range("I102").Value = 10000
$ws.Range("K102").Value = 10000
$ws.Range("M102").Value = -6755
# row 105
$ws.Range("H105").Value = 2482.9565
$ws.Range("I105").Value = 2459.4546
$ws.Range("K105").Value = 2459.4546
$ws.Range("M105").Value = -712.4546
# row 134
$ws.Range("H134").Value = 6735.394
$ws.Range("I134").Value = 6737.643
$ws.Range("K134").Value = 20212.929
$ws.Range("M134").Value = -17677.929

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
# row 31
$ws.Range("H31").Value = 2363.0557
$ws.Range("I31").Value = 2117.5715
$ws.Range("K31").Value = 2117.5715
$ws.Range("M31").Value = -1822.5715
# row 34
$ws.Range("H34").Value = 2363.0557
$ws.Range("I34").Value = 2117.5715
$ws.Range("K34").Value = 2117.5715
$ws.Range("M34").Value = -1915.5715

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
# row 70
$ws.Range("H70").Value = 42150
$ws.Range("J70").Value = 4300
$ws.Range("L70").Value = 4300
$ws.Range("N70").Value = -4840
# row 73
$ws.Range("H73").Value = 42150
$ws.Range("J73").Value = 4300
$ws.Range("L73").Value = 4300
$ws.Range("N73").Value = -6172
# row 80
$ws.Range("H80").Value = 3150.2222
$ws.Range("J80").Value = 3725.5
$ws.Range("L80").Value = 3725.5
$ws.Range("N80").Value = -5721.5
# row 83
$ws.Range("H83").Value = 3150.2222
$ws.Range("J83").Value = 3725.5
$ws.Range("L83").Value = 18627.5
$ws.Range("N83").Value = -28611.5
# row 98
$ws.Range("H98").Value = 31830
$ws.Range("J98").Value = 31830
$ws.Range("L98").Value = 31830
$ws.Range("N98").Value = -37820
# row 126
$ws.Range("H126").Value = 2177782.8
$ws.Range("J126").Value = 57946.89
$ws.Range("L126").Value = 173840.67
$ws.Range("N126").Value = -178780.67
# row 132
$ws.Range("H132").Value = 1041970.06
$ws.Range("I132").Value = 1604202
$ws.Range("J132").Value = 4003.4614
$ws.Range("K132").Value = 4812606
$ws.Range("L132").Value = 12010.3842
$ws.Range("M132").Value = -4810076
$ws.Range("N132").Value = -17070.3842
# row 134
$ws.Range("H134").Value = 23941.111
$ws.Range("J134").Value = 23941.111
$ws.Range("L134").Value = 71823.333
$ws.Range("N134").Value = -76893.333
# row 136
$ws.Range("H136").Value = 9985.666999999999
$ws.Range("J136").Value = 9985.666999999999
$ws.Range("L136").Value = 29957.001
$ws.Range("N136").Value = -35057.001

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
# row 82
$ws.Range("H82").Value = 1940
$ws.Range("I82").Value = 1442.5
$ws.Range("K82").Value = 1442.5
$ws.Range("M82").Value = -1081.5
# row 85
$ws.Range("H85").Value = 1940
$ws.Range("I85").Value = 1442.5
$ws.Range("K85").Value = 1442.5
$ws.Range("M85").Value = -194.5
# row 101
$ws.Range("H101").Value = 6499.5
$ws.Range("J101").Value = 6499.5
$ws.Range("L101").Value = 6499.5
$ws.Range("N101").Value = -12989.5
# row 122
$ws.Range("H122").Value = 4752.533
$ws.Range("J122").Value = 3975.375
$ws.Range("L122").Value = 11926.125
$ws.Range("N122").Value = -16826.125

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
# row 63
$ws.Range("H63").Value = 27999.334
$ws.Range("J63").Value = 27999.334
$ws.Range("L63").Value = 27999.334
$ws.Range("N63").Value = -29247.334
# row 66
$ws.Range("H66").Value = 27999.334
$ws.Range("J66").Value = 27999.334
$ws.Range("L66").Value = 83998.00199999999
$ws.Range("N66").Value = -90238.00199999999
# row 113
$ws.Range("H113").Value = 578.3889
$ws.Range("I113").Value = 279.42856
$ws.Range("K113").Value = 838.28568
$ws.Range("M113").Value = 1331.71432
# row 132
$ws.Range("H132").Value = 1200.7747
$ws.Range("I132").Value = 837.4032
$ws.Range("J132").Value = 3704
$ws.Range("K132").Value = 2512.2096
$ws.Range("L132").Value = 11112
$ws.Range("M132").Value = 17.79039999999986
$ws.Range("N132").Value = -16172
